$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("SCHMITT", "Hadrien", "21 Bd Maréchal Lyautey, 38000 Grenoble", 45.1859198, 5.731540139877457),
    @("SCHMITT", "Hadrien", "21 Bd Maréchal Lyautey, 38000 Grenoble", 45.1859198, 5.731540139877457),
    @("ASTRID", "Monique", "30 Rue Rambaud, 17000 La Rochelle", 46.16360619065998, -1.155014376554859)
)

$row = 7
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
